$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1825396825396825
$ws.Range("C2").Value = 0.5912698412698413
$ws.Range("J2").Value = 0.0119047619047619
$ws.Range("P2").Value = 0.126984126984127
$ws.Range("S2").Value = 0.0873015873015873
$ws.Range("B3").Value = 0.02531645569620253
$ws.Range("C3").Value = 0.04430379746835443
$ws.Range("J3").Value = 0.0379746835443038
$ws.Range("P3").Value = 0.689873417721519
$ws.Range("S3").Value = 0.2025316455696203
$ws.Range("J4").Value = 0.02380952380952381
$ws.Range("P4").Value = 0.7380952380952381
$ws.Range("S4").Value = 0.2380952380952381
$ws.Range("B6").Value = 0.04784688995215311
$ws.Range("D6").Value = 0.02392344497607655
$ws.Range("F6").Value = 0.03827751196172249
$ws.Range("J6").Value = 0.2296650717703349
$ws.Range("O6").Value = 0.01913875598086124
$ws.Range("Q6").Value = 0.1770334928229665
$ws.Range("R6").Value = 0.1004784688995215
$ws.Range("S6").Value = 0.3636363636363636
$ws.Range("B7").Value = 0.04615384615384616
$ws.Range("D7").Value = 0.003846153846153846
$ws.Range("E7").Value = 0.003846153846153846
$ws.Range("F7").Value = 0.04230769230769231
$ws.Range("J7").Value = 0.1384615384615385
$ws.Range("O7").Value = 0.01538461538461539
$ws.Range("Q7").Value = 0.1653846153846154
$ws.Range("R7").Value = 0.08846153846153847
$ws.Range("S7").Value = 0.4961538461538462
$ws.Range("B8").Value = 0.0562015503875969
$ws.Range("D8").Value = 0.01550387596899225
$ws.Range("F8").Value = 0.03875968992248062
$ws.Range("J8").Value = 0.1143410852713178
$ws.Range("O8").Value = 0.02131782945736434
$ws.Range("Q8").Value = 0.1744186046511628
$ws.Range("R8").Value = 0.1220930232558139
$ws.Range("S8").Value = 0.4573643410852713
$ws.Range("B9").Value = 0.08203125
$ws.Range("D9").Value = 0.0078125
$ws.Range("F9").Value = 0.046875
$ws.Range("J9").Value = 0.10546875
$ws.Range("O9").Value = 0.015625
$ws.Range("Q9").Value = 0.19921875
$ws.Range("R9").Value = 0.109375
$ws.Range("S9").Value = 0.43359375
$ws.Range("B10").Value = 0.08852005532503458
$ws.Range("D10").Value = 0.01798063623789765
$ws.Range("E10").Value = 0.0006915629322268327
$ws.Range("F10").Value = 0.06500691562932227
$ws.Range("J10").Value = 0.1154910096818811
$ws.Range("O10").Value = 0.01452282157676349
$ws.Range("Q10").Value = 0.210926694329184
$ws.Range("R10").Value = 0.08644536652835408
$ws.Range("S10").Value = 0.4004149377593361
$ws.Range("G11").Value = 0.1359773371104816
$ws.Range("J11").Value = 0.06515580736543909
$ws.Range("K11").Value = 0.1954674220963173
$ws.Range("L11").Value = 0.5864022662889519
$ws.Range("S11").Value = 0.0169971671388102
$ws.Range("G12").Value = 0.8425925925925926
$ws.Range("J12").Value = 0.09722222222222222
$ws.Range("K12").Value = 0.004629629629629629
$ws.Range("L12").Value = 0.03240740740740741
$ws.Range("S12").Value = 0.02314814814814815
$ws.Range("G13").Value = 0.6851851851851852
$ws.Range("J13").Value = 0.2777777777777778
$ws.Range("S13").Value = 0.03703703703703703
$ws.Range("F15").Value = 0.01310043668122271
$ws.Range("H15").Value = 0.1572052401746725
$ws.Range("I15").Value = 0.07423580786026202
$ws.Range("J15").Value = 0.3930131004366812
$ws.Range("K15").Value = 0.05676855895196507
$ws.Range("M15").Value = 0.01746724890829694
$ws.Range("O15").Value = 0.03930131004366812
$ws.Range("S15").Value = 0.2489082969432314
$ws.Range("F16").Value = 0.02409638554216868
$ws.Range("H16").Value = 0.1807228915662651
$ws.Range("I16").Value = 0.0963855421686747
$ws.Range("J16").Value = 0.4216867469879518
$ws.Range("K16").Value = 0.1265060240963855
$ws.Range("M16").Value = 0.03614457831325301
$ws.Range("O16").Value = 0.02409638554216868
$ws.Range("S16").Value = 0.09036144578313253
$ws.Range("F17").Value = 0.02095238095238095
$ws.Range("H17").Value = 0.1733333333333333
$ws.Range("I17").Value = 0.1123809523809524
$ws.Range("J17").Value = 0.4228571428571429
$ws.Range("K17").Value = 0.09333333333333334
$ws.Range("M17").Value = 0.02857142857142857
$ws.Range("O17").Value = 0.05523809523809524
$ws.Range("S17").Value = 0.09333333333333334
$ws.Range("F18").Value = 0.01532567049808429
$ws.Range("H18").Value = 0.1800766283524904
$ws.Range("I18").Value = 0.09195402298850575
$ws.Range("J18").Value = 0.4137931034482759
$ws.Range("K18").Value = 0.08812260536398467
$ws.Range("M18").Value = 0.01149425287356322
$ws.Range("N18").Value = 0.003831417624521073
$ws.Range("O18").Value = 0.07279693486590039
$ws.Range("S18").Value = 0.1226053639846743
$ws.Range("F19").Value = 0.01049868766404199
$ws.Range("H19").Value = 0.2106299212598425
$ws.Range("I19").Value = 0.09448818897637795
$ws.Range("J19").Value = 0.3700787401574803
$ws.Range("K19").Value = 0.1161417322834646
$ws.Range("M19").Value = 0.01837270341207349
$ws.Range("N19").Value = 0.001312335958005249
$ws.Range("O19").Value = 0.06233595800524935
$ws.Range("S19").Value = 0.1161417322834646
